$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2-26 down to 3-27.
$ws.Range("A2:R2").Insert() | Out-Null

# The inserted range picks up the header row's formatting by default; reset it to
# plain/unstyled like the rest of the data rows, then re-apply the date format used
# in column D for all data rows.
$ws.Range("A2:R2").Style = "Normal"
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Fill in the values for new row 2 (copy of old row 2 / now row 3, but with date 44631).
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(2, 3).Value = "Maule"
$ws.Cells.Item(2, 4).Value = 44631
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 100112043
$ws.Cells.Item(2, 7).Value = "Pepino dulce"
$ws.Cells.Item(2, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 15000
$ws.Cells.Item(2, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 833
$ws.Cells.Item(2, 17).Value = 18
$ws.Cells.Item(2, 18).Value = "Hortaliza"
